$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume) hold plain text that looks numeric
# (e.g. "1.00", "0.0944", "  -2.53%  "). A direct .Value assignment lets
# Excel auto-convert these to real numbers/percentages, which would lose
# the exact original text (trailing zeros, spacing, % sign). Force the
# cell to Text first, write the literal string, then restore the original
# (default) style so no stray formatting is left behind.
function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

# Row 2
Set-TextValue "D2" "54.996.08"
Set-TextValue "E2" "  -3.64%  "

# Row 3
Set-TextValue "D3" "2.330.51"
Set-TextValue "E3" "  -2.78%  "

# Row 4
Set-TextValue "E4" "  +0.15%  "

# Row 5
Set-TextValue "D5" "496.46"
Set-TextValue "E5" "  -2.38%  "

# Row 6
Set-TextValue "D6" "129.38"
Set-TextValue "E6" "  -3.90%  "

# Row 7
Set-TextValue "E7" "  +0.22%  "

# Row 8
Set-TextValue "D8" "0.530"
Set-TextValue "E8" "  -4.45%  "

# Row 9
Set-TextValue "D9" "2.339.17"
Set-TextValue "E9" "  -2.73%  "

# Row 10
Set-TextValue "D10" "0.0944"
Set-TextValue "E10" "  -4.43%  "

# Row 11
Set-TextValue "E11" "  -1.32%  "

# Row 12
Set-TextValue "D12" "4.73"
Set-TextValue "E12" "  +1.08%  "

# Row 13
Set-TextValue "D13" "0.316"
Set-TextValue "E13" "  -5.41%  "

# Row 14
Set-TextValue "D14" "2.749.10"
Set-TextValue "E14" "  -2.66%  "

# Row 15
$ws.Range("B15").Value = "Avalanche"
$ws.Range("C15").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextValue "D15" "21.39"
Set-TextValue "E15" "  -2.43%  "

# Row 16
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
Set-TextValue "D16" "54.978.87"
Set-TextValue "E16" "  -3.53%  "

# Row 17
Set-TextValue "D17" "0.0000129"
Set-TextValue "E17" "  -3.49%  "

# Row 18
Set-TextValue "D18" "2.329.43"
Set-TextValue "E18" "  -2.18%  "

# Row 19
Set-TextValue "D19" "9.70"
Set-TextValue "E19" "  -4.87%  "

# Row 20
Set-TextValue "D20" "305.22"
Set-TextValue "E20" "  -2.07%  "

# Row 21
Set-TextValue "D21" "3.97"
Set-TextValue "E21" "  -2.15%  "

# Row 22
Set-TextValue "E22" "  -0.98%  "

# Row 23
Set-TextValue "E23" "  +0.04%  "

# Row 24
Set-TextValue "D24" "64.40"
Set-TextValue "E24" "  -1.31%  "

# Row 25
Set-TextValue "D25" "0.998"
Set-TextValue "E25" "  +0.32%  "

# Row 26
Set-TextValue "D26" "0.368"
Set-TextValue "E26" "  -1.96%  "

# Row 27
Set-TextValue "D27" "0.144"
Set-TextValue "E27" "  -5.05%  "

# Row 28
Set-TextValue "D28" "7.11"
Set-TextValue "E28" "  -3.99%  "

# Row 29
Set-TextValue "D29" "167.93"
Set-TextValue "E29" "  -2.52%  "

# Row 30
Set-TextValue "D30" "0.0₃0700"
Set-TextValue "E30" "  -4.51%  "

# Row 31
Set-TextValue "D31" "1.62"
Set-TextValue "E31" "  -1.92%  "

# Row 32
Set-TextValue "D32" "0.999"
Set-TextValue "E32" "  -0.06%  "

# Row 33
Set-TextValue "D33" "1.00"
Set-TextValue "E33" "  +0.34%  "

# Row 34
Set-TextValue "D34" "5.70"
Set-TextValue "E34" "  -2.96%  "

# Row 35
Set-TextValue "D35" "1.06"
Set-TextValue "E35" "  -6.05%  "

# Row 36
Set-TextValue "D36" "17.53"
Set-TextValue "E36" "  -2.13%  "

# Row 37
Set-TextValue "E37" "  -3.72%  "

# Row 38
Set-TextValue "D38" "0.843"
Set-TextValue "E38" "  +2.28%  "

# Row 39
Set-TextValue "D39" "3.62"
Set-TextValue "E39" "  -6.71%  "

# Row 40
Set-TextValue "D40" "35.97"
Set-TextValue "E40" "  -0.59%  "

# Row 41
Set-TextValue "D41" "0.370"
Set-TextValue "E41" "  -1.23%  "

# Row 42
Set-TextValue "D42" "1.37"
Set-TextValue "E42" "  -4.89%  "

# Row 43
Set-TextValue "D43" "3.33"
Set-TextValue "E43" "  -2.27%  "

# Row 44
$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue "D44" "123.71"
Set-TextValue "E44" "  -6.59%  "

# Row 45
$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D45" "4.72"
Set-TextValue "E45" "  -3.35%  "

# Row 46
Set-TextValue "D46" "0.549"
Set-TextValue "E46" "  -3.58%  "

# Row 47
Set-TextValue "D47" "0.0885"
Set-TextValue "E47" "  -3.01%  "

# Row 48
Set-TextValue "D48" "239.03"
Set-TextValue "E48" "  -4.39%  "

# Row 49
Set-TextValue "D49" "0.0475"
Set-TextValue "E49" "  -2.67%  "

# Row 50
Set-TextValue "D50" "16.67"
Set-TextValue "E50" "  -2.44%  "

# Row 51
Set-TextValue "E51" "  -2.84%  "
